$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 183.76471
$ws.Range("I53").Value = 93.333336
$ws.Range("J53").Value = 285.5
$ws.Range("K53").Value = 93.333336
$ws.Range("L53").Value = 285.5
$ws.Range("M53").Value = 543.666664
$ws.Range("N53").Value = -1559.5
$ws.Range("H98").Value = 11571.363
$ws.Range("I98").Value = 9140.714
$ws.Range("J98").Value = 15825
$ws.Range("K98").Value = 9140.714
$ws.Range("L98").Value = 15825
$ws.Range("M98").Value = -7642.714
$ws.Range("N98").Value = -18821
$ws.Range("H116").Value = 13335580
$ws.Range("I116").Value = 66668300
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 66668300
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = -66664858
$ws.Range("N116").Value = -9284
$ws.Range("H122").Value = 11571.363
$ws.Range("I122").Value = 9140.714
$ws.Range("J122").Value = 15825
$ws.Range("K122").Value = 27422.142
$ws.Range("L122").Value = 47475
$ws.Range("M122").Value = -24972.142
$ws.Range("N122").Value = -52375
$ws.Range("H125").Value = 4291.143
$ws.Range("I125").Value = 4766.3335
$ws.Range("K125").Value = 42897.0015
$ws.Range("M125").Value = -40437.0015
$ws.Range("H129").Value = 1103.9104
$ws.Range("J129").Value = 1164.7097
$ws.Range("L129").Value = 3494.1291
$ws.Range("N129").Value = -13494.1291
$ws.Range("H138").Value = 2065.054
$ws.Range("I138").Value = 1381.8667
$ws.Range("J138").Value = 4993
$ws.Range("K138").Value = 4145.6001
$ws.Range("L138").Value = 14979
$ws.Range("M138").Value = 994.3999000000003
$ws.Range("N138").Value = -25259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 11333.333
$ws.Range("I6").Value = 11333.333
$ws.Range("K6").Value = 11333.333
$ws.Range("M6").Value = -11160.333
$ws.Range("H23").Value = 10907.692
$ws.Range("I23").Value = 9988.888999999999
$ws.Range("K23").Value = 9988.888999999999
$ws.Range("M23").Value = -9729.888999999999
$ws.Range("H32").Value = 14049.036
$ws.Range("I32").Value = 15442.6875
$ws.Range("J32").Value = 4492.5713
$ws.Range("K32").Value = 15442.6875
$ws.Range("L32").Value = 4492.5713
$ws.Range("M32").Value = -15155.6875
$ws.Range("N32").Value = -5066.5713
$ws.Range("H44").Value = 39000
$ws.Range("J44").Value = 39000
$ws.Range("L44").Value = 39000
$ws.Range("N44").Value = -39976
$ws.Range("H55").Value = 28333.334
$ws.Range("H80").Value = 34925
$ws.Range("J80").Value = 34925
$ws.Range("L80").Value = 34925
$ws.Range("N80").Value = -36921
$ws.Range("H83").Value = 34925
$ws.Range("J83").Value = 34925
$ws.Range("L83").Value = 104775
$ws.Range("N83").Value = -114759

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2540.6667
$ws.Range("I94").Value = 2222
$ws.Range("J94").Value = 2700
$ws.Range("K94").Value = 2222
$ws.Range("L94").Value = 2700
$ws.Range("M94").Value = -1771
$ws.Range("N94").Value = -3602
$ws.Range("H134").Value = 2030.2084
$ws.Range("I134").Value = 2036.3
$ws.Range("J134").Value = 1999.75
$ws.Range("K134").Value = 6108.9
$ws.Range("L134").Value = 5999.25
$ws.Range("M134").Value = -3573.9
$ws.Range("N134").Value = -11069.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 378.85715
$ws.Range("I22").Value = 231.33333
$ws.Range("K22").Value = 231.33333
$ws.Range("M22").Value = 118.66667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2440.6
$ws.Range("I33").Value = 700.3333
$ws.Range("J33").Value = 5051
$ws.Range("K33").Value = 4201.9998
$ws.Range("L33").Value = 30306
$ws.Range("M33").Value = -3918.9998
$ws.Range("N33").Value = -30872
$ws.Range("H131").Value = 20456.424
$ws.Range("I131").Value = 362.30768
$ws.Range("J131").Value = 27154.46
$ws.Range("K131").Value = 1086.92304
$ws.Range("L131").Value = 81463.38
$ws.Range("M131").Value = 3953.07696
$ws.Range("N131").Value = -91543.38
$ws.Range("H133").Value = 4576.95
$ws.Range("I133").Value = 1732.625
$ws.Range("J133").Value = 6473.1665
$ws.Range("K133").Value = 5197.875
$ws.Range("L133").Value = 19419.4995
$ws.Range("M133").Value = -137.875
$ws.Range("N133").Value = -29539.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33622
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -104112
$ws.Range("H132").Value = 3936
$ws.Range("I132").Value = 3417.8
$ws.Range("J132").Value = 5231.5
$ws.Range("K132").Value = 10253.4
$ws.Range("L132").Value = 15694.5
$ws.Range("M132").Value = -7723.400000000001
$ws.Range("N132").Value = -20754.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4812.875
$ws.Range("I7").Value = 4669.3335
$ws.Range("J7").Value = 4899
$ws.Range("K7").Value = 4669.3335
$ws.Range("L7").Value = 4899
$ws.Range("M7").Value = -4557.3335
$ws.Range("N7").Value = -5123
$ws.Range("H22").Value = 596.4286
$ws.Range("I22").Value = 515
$ws.Range("K22").Value = 515
$ws.Range("M22").Value = -220
$ws.Range("H27").Value = 596.4286
$ws.Range("I27").Value = 515
$ws.Range("K27").Value = 515
$ws.Range("M27").Value = -408
$ws.Range("H76").Value = 6324.5
$ws.Range("J76").Value = 7432.6665
$ws.Range("L76").Value = 7432.6665
$ws.Range("N76").Value = -8108.6665
$ws.Range("H79").Value = 6324.5
$ws.Range("J79").Value = 7432.6665
$ws.Range("L79").Value = 7432.6665
$ws.Range("N79").Value = -9772.666499999999
$ws.Range("H126").Value = 4812.875
$ws.Range("I126").Value = 4669.3335
$ws.Range("J126").Value = 4899
$ws.Range("K126").Value = 14008.0005
$ws.Range("L126").Value = 14697
$ws.Range("M126").Value = -11538.0005
$ws.Range("N126").Value = -19637
$ws.Range("H132").Value = 3568.2368
$ws.Range("I132").Value = 3359.8572
$ws.Range("K132").Value = 10079.5716
$ws.Range("M132").Value = -7549.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 32373.04
$ws.Range("J123").Value = 32373.04
$ws.Range("L123").Value = 32373.04
$ws.Range("N123").Value = -42173.04
